$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 82, shifting existing rows 82-183 down to 83-184
$ws.Rows(82).Insert()

# Populate the newly inserted row 82 with the new record
$ws.Cells.Item(82, 1).Value = 5
$ws.Cells.Item(82, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value = "Maule"
$ws.Range("D82").Value = 44494
$ws.Cells.Item(82, 5).Value = 7
$ws.Cells.Item(82, 6).Value = 100112006
$ws.Cells.Item(82, 7).Value = "Repollo"
$ws.Cells.Item(82, 8).Value = "Crespo record"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 5000
$ws.Cells.Item(82, 11).Value = 600
$ws.Cells.Item(82, 12).Value = 600
$ws.Cells.Item(82, 13).Value = 600
$ws.Cells.Item(82, 14).Value = "`$/unidad"
$ws.Cells.Item(82, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(82, 16).Value = 600
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"
